$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the RF value (column I) for rows 31 through 66 from the old
# 2025 reference factor to the newly recalculated one.
$newRF = 10.07971014492754

for ($r = 31; $r -le 66; $r++) {
    $ws.Cells.Item($r, 9).Value = $newRF
}
